$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the contents of rows 2-19 (keeps existing cell formatting, e.g.
# the bold/bordered style on column A), and fully clear rows 20-32
# (formatting included) so that once empty they disappear from the
# sheet entirely and the used range shrinks down to A1:F19.
$ws.Range("A2:F19").ClearContents()
$ws.Range("A20:F32").Clear()

# New data for rows 2-19 (columns B=Buying Opportunity, C=support Zone,
# D=long buildup, E=Short buildup, F=FII ENTERING). Column A keeps the
# original 0-based index numbers.
$data = @(
    @(0,  "NSE:3IINFOLTD",  "NSE:BFUTILITIE", "", "", "NSE:GODREJPROP"),
    @(1,  "NSE:AROGRANITE", "NSE:HPAL",       "", "", "NSE:KFINTECH"),
    @(2,  "NSE:EVERESTIND", "NSE:JKPAPER",    "", "", "NSE:OBEROIRLTY"),
    @(3,  "NSE:FDC",        "NSE:MBLINFRA",   "", "", ""),
    @(4,  "NSE:GODREJPROP", "NSE:MUNJALAU",   "", "", ""),
    @(5,  "NSE:HDFCLIQUID", "NSE:PALASHSECU", "", "", ""),
    @(6,  "NSE:HEIDELBERG", "NSE:PENINLAND",  "", "", ""),
    @(7,  "NSE:ICDSLTD",    "",               "", "", ""),
    @(8,  "NSE:INTENTECH",  "",               "", "", ""),
    @(9,  "NSE:IOLCP",      "",               "", "", ""),
    @(10, "NSE:JSL",        "",               "", "", ""),
    @(11, "NSE:KICL",       "",               "", "", ""),
    @(12, "NSE:LOVABLE",    "",               "", "", ""),
    @(13, "NSE:MANKIND",    "",               "", "", ""),
    @(14, "NSE:MONIFTY500", "",               "", "", ""),
    @(15, "NSE:NETWEB",     "",               "", "", ""),
    @(16, "NSE:ORCHPHARMA", "",               "", "", ""),
    @(17, "NSE:PNBHOUSING", "",               "", "", "")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $ws.Cells.Item($row, 6).Value = $entry[5]
    $row++
}
